$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 02:59"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8215199
$ws.Range("C4").Value = 65013
$ws.Range("D4").Value = 5318012
$ws.Range("E4").Value = 2674475
$ws.Range("G4").Value = 869
$ws.Range("H4").Value = 222712

# Row 11 - Peru
$ws.Range("B11").Value = 859740
$ws.Range("C11").Value = 2789
$ws.Range("D11").Value = 764750
$ws.Range("E11").Value = 61413
$ws.Range("G11").Value = 65
$ws.Range("H11").Value = 33577

# Row 30 - Paises Bajos
$ws.Range("B30").Value = 191732
$ws.Range("C30").Value = 2345
$ws.Range("E30").Value = 20543

# Row 39 - Bolivia
$ws.Range("B39").Value = 122883
$ws.Range("C39").Value = 755
$ws.Range("D39").Value = 98671
$ws.Range("E39").Value = 21683
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 2529

# Row 69 - Armenia
$ws.Range("B69").Value = 47173
$ws.Range("C69").Value = 47
$ws.Range("D69").Value = 46527
$ws.Range("E69").Value = 336

# Row 117
$ws.Range("B117").Value = 7585
$ws.Range("C117").Value = 13
$ws.Range("D117").Value = 7326
$ws.Range("E117").Value = 96

# Row 118
$ws.Range("B118").Value = 7444
$ws.Range("C118").Value = 73
$ws.Range("D118").Value = 6348
$ws.Range("E118").Value = 1017

# Row 147
$ws.Range("B147").Value = 3620
$ws.Range("C147").Value = 31
$ws.Range("D147").Value = 2522
$ws.Range("E147").Value = 991
$ws.Range("G147").Value = 1
$ws.Range("H147").Value = 107

# Row 154
$ws.Range("B154").Value = 2496
$ws.Range("C154").Value = 18
$ws.Range("E154").Value = 125

# Row 183
$ws.Range("B183").Value = 422
$ws.Range("C183").Value = 8
$ws.Range("D183").Value = 376
$ws.Range("E183").Value = 46

# Row 184
$ws.Range("B184").Value = 415
$ws.Range("C184").Value = 8
$ws.Range("D184").Value = 364
$ws.Range("E184").Value = 41

# Row 191
$ws.Range("B191").Value = 218
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 195
